$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text (A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.32 = 12776.5 pesos`n✅ 12776.5 pesos = 3.29 = 929.1 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 301.57
$ws2.Range("O10").Value = 3853.01
$ws2.Range("N12").Value = 3880
$ws2.Range("O12").Value = 282.151
